# Applies Xhosa translation edits per the commit diff.
$d = $word.ActiveDocument

function Replace-All($range, $old, $new) {
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Simple whole-run text replacements across the document body ---
Replace-All $d.Content "Appendix 14: SWIFT Consent Call to Parent for Adolescent Interview " "ISihlomelo 14: Umnxeba weMvume oya kuMzali kudliwano-ndlebe lwaBafikisayo lwe-SWIFT "
Replace-All $d.Content "Before making this call a member of the team should have sent a copy of the consent form to the parent along with the following message from the study WhatsApp number. " "Ngaphambi kokuba utsalele lomnxeba ilungu leqela bekumele ukuba lithumele ikopi yefomu yemvume kumzali kunye nalomyalezo ulandelayo osuka kwinombolo yophononongo kaWhatsApp. "
Replace-All $d.Content "Hello. This is _______ from the SWIFT research team looking at the ParentText Chatbot you have been interacting with. I have sent you a copy of the information form which you got when you said that we could interview your child a while back. They have been selected for the second interview which we explained in the form. One of our team will be calling you to go over the information, get your consent and set up a time to call your child. Please have a look at it before the call, if possible, so that you can ask any questions you have. Thanks!" "Molo. This is _______ from the SWIFT research team looking at the ParentText Chatbot you have been interacting with. Ndikuthumelele ikopi yefomu yolwazi owayifumanayo xa ubusithi singenza udliwano-ndlebe nomntwana wakho kwithutyana elidlulileyo. Uye wakhethwa kudliwano-ndlebe lwesibini esiluchaze kwifomu. Elinye leqela lethu liya kukutsalela umnxeba ukuba ujonge ulwazi, sifumane imvume yakho kwaye umisele ixesha lokutsalela umntwana wakho. Nceda uyijonge ngaphambi komnxeba, ukuba kuyenzeka, ukuze ukwazi ukubuza nayiphi na imibuzo onayo. Enkosi!"
Replace-All $d.Content "Call script  " "Isikripthi somnxeba  "
Replace-All $d.Content "Hi there.  I’m ____________. I’m calling from the SWIFT research team looking at the ParentText Chatbot you have been interacting with. A while back, you gave consent for your child to be interviewed by our team. As you may know we have already done the shorter survey interview, but your child has also been selected for the one-on-one telephonic interview. You would have received a message from our team with the information sheet, saying that we will be calling you. Is this a good time for you?" "Molo apho.  Ndingu ____________. I’m calling from the SWIFT research team looking at the ParentText Chatbot you have been interacting with. Kwithutyana elidlulileyo, uye wasinika imvume yokuba umntwana wakho enziwe udliwano-ndlebe liqela lethu. Njengoko usazi, sele silwenzile udliwano-ndlebe olufutshane lwesaveyi, kodwa umntwana wakho naye ukwakhethelwe udliwano-ndlebe lwabucala ngomnxeba. Kumele uwufumene umyalezo ovela kwiqela lethu kunye nephepha lolwazi, lusithi sizakufowunela. Ingaba lixesha elilungileyo eli kuwe?"
Replace-All $d.Content "If they answer no, get an alternate time to call back, if yes proceed." "Ukuba uphendule ngo hayi, fumana elinye ixesha lokufowuna, ukuba uthi ewe qhubeka."
Replace-All $d.Content "I’m calling to get your consent to interview ___________ and before you give it to me, I wanted to make sure you have gone through the consent form to make sure you know all the details you need to know to give that consent. Have you read it already or would you like me to go through it with you now? " "Ndikutsalela umnxeba ukufumana imvume yakho yokwenza udliwano-ndlebe no___________ kwaye phambi kokuba undinike, bendifuna ukuqinisekisa ukuba uyijongile ifomu yemvume ukuqiniseka uyazazi zonke iinkcukacha odinga ukuzazi ukunika loo mvume. Ingaba sele uyifundile okanye ungathanda ukuba ndiyijonge nawe ngoku? "
Replace-All $d.Content "Please stop me as I’m going if anything isn’t clear or if you have any questions. I will also make sure that all your questions are answered at the end. " "Nceda undimise njengoko ndiqhuba ukuba kukho into engacacanga okanye ukuba unemibuzo. Ndiza kuqinisekisa ukuba yonke imibuzo yakho iyaphendulwa ekugqibeleni. "
Replace-All $d.Content "Proceed to read consent form slowly" "Qhubeka nokufunda ifomu yemvume ngokucothayo"
Replace-All $d.Content "I’m glad that you went through it. Do you have any questions? " "Ndiyavuya ukuba uye wayijonga. Unayo nayiphi na imibuzo? "
Replace-All $d.Content "Do you have any questions?" "Unayo nayiphi na imibuzo?"
Replace-All $d.Content "Proceed to informed consent to take part in study below" "Qhubeka kwimvume enolwazi yokuthatha inxaxheba kuphononongo olungezantsi"
Replace-All $d.Content "I know that you have read it, but I would just like to go over the last part which is the consent part. " "Ndiyayazi ukuba uyifundile, kodwa ndingathanda ukujonga kwakhona kwindawo yokugqibela eyinxalenye yemvume. "
Replace-All $d.Content "Proceed to informed consent to take part in study below" "Qhubeka kwimvume enolwazi yokuthatha inxaxheba kuphononongo olungezantsi"
Replace-All $d.Content "Informed Consent to Take Part in the Study" "Imvume eChaziweyo yokuThatha Inxaxheba kuPhononongo"
Replace-All $d.Content "Please listen carefully to the following questions and reply yes if you consent. If you say no to anything, that's okay. I will go over that bit of information again and then you can confirm once you are happy." "Nceda umamele ngononophelo kulemibuzo ilandelayo kwaye uphendule ngo ewe ukuba uyavuma. Ukuba uthi hayi nakwintoni na, oko kulungile. Ndiza kuphinda ndijonge kancinci ulwazi kwakhona kwaye emvakoko unokuqinisekisa xa wonwabile."
Replace-All $d.Content "Do you confirm that you have read the information in the consent form and know what is expected of your child?" "Ingaba uyaqinisekisa ukuba ulufundile ulwazi olukwifomu yemvume kwaye uyakwazi okulindeleke kumntwana wakho?"
Replace-All $d.Content "Do you understand as your child’s guardian that you are giving consent for them to participate? " "Ingaba uyayiqonda ukuba njengomnonopheli womntwana wakho unikeza imvume yokuba athathe inxaxheba? "
Replace-All $d.Content "Do you understand that even though you have given consent, that your child will still be able to choose freely if they want to be interviewed? " "Ingaba uyayiqonda into yokuba nangona uyinikezile imvume, umntwana wakho useza kukwazi ukukhetha ngokukhululekileyo ukuba uyalufuna udliwano-ndlebe? "
Replace-All $d.Content "Do you understand that they can say no to being interviewed without any consequence? " "Ingaba uyayiqonda into yokuba banokuthi hayi kudliwano-ndlebe ngaphandle kwaso nasiphi na isiphumo? "
Replace-All $d.Content "Have you had time to think about the information?" "Ingaba uye walifumana ixesha lokucinga malunga nolwazi?"
Replace-All $d.Content "Have you asked any questions you had?" "Ingaba uye wayibuza imibuzo ubunayo?"
Replace-All $d.Content "Are you satisfied with the answers you got if you did have questions?" "Ingaba wanelisekile ziimpendulo ozifumeneyo ukuba ubunayo imibuzo?"
Replace-All $d.Content "Do you know who can see your child’s information, how it will be kept safe, and what happens to it after the study?" "Ingaba uyayazi ukuba ngubani onokubona ulwazi lomntwana wakho, ukuba luya kugcinwa njani lukhuselekile, kwaye kwenzeka ntoni kulo emva kophononongo?"
Replace-All $d.Content "Do you understand that you will not be notified of your child’s answers?" "Ingaba uyayiqonda ukuba awuyi kwaziswa ngeempendulo zomntwana wakho?"
Replace-All $d.Content "Do you understand that if the researchers pick up any safety concerns while interviewing your child, they will let your child know that they will have to share the information, and then they will share it with you?" "Ingaba uyayiqonda into yokuba ukuba abaphandi bafumanisa naziphi na iinkxalabo zokhuseleko ngelixa lodliwano-ndlebe nomntwana wakho, baya kumazisa umntwana wakho ukuba kuya kufuneka babelane ngolwazi, kwaye emvakoko baya kwabelana nawe ngalo?"
Replace-All $d.Content "Do you understand that you and your child won’t be named in any results, papers or reports from this study?" "Ingaba uyayiqonda ukuba wena nomntwana wakho anizokuchazwa kuzo naziphi na iziphumo, amaphepha okanye iingxelo zolu phononongo?"
Replace-All $d.Content "Do you know who to contact if you have a problem with the study?" "Ingaba uyayazi ukuba ungaqhagamshelana nabani ukuba unengxaki ngoluphononongo?"
Replace-All $d.Content "Can one of the team contact you again if more information is needed from you?" "Ingaba omnye weqela angaqhagamshelana nawe kwakhona ukuba ulwazi oluninzi luyafuneka kuwe?"
Replace-All $d.Content "Can we keep your contact information so we can tell you about the results of the study?" "Ingaba singazigcina iinkcukacha zakho zoqhagamshelwano ukuze sikuxelele ngeziphumo zophando?"
Replace-All $d.Content "Do you give us consent to contact your child and for them to take part in the interview? " "Ingaba uyasinika imvume yokuba siqhagamshelane nomntwana wakho kwaye athathe inxaxheba kudliwano-ndlebe? "
Replace-All $d.Content "I would like to set up a time which will work for us to call your child for the interview. But I just want to make sure that when we do call, they will be able to speak privately in a space where they won’t be disturbed or overheard. Would you please be able to help us make sure that they are able to talk without anyone overhearing them, even you. This is to ensure their privacy and confidentiality. I know it’s difficult in some of our houses, but can you think of a time and space where this is possible? (explore this space with them or how to perhaps put something in place)" "Ndingathanda ukucwangcisa ixesha eliya kuthi lisisebenzele ukuze sitsalele umnxeba umntwana wakho kudliwano-ndlebe. Kodwa ndifuna nje ukuqinisekisa ukuba xa simtsalela umnxeba, uya kukwazi ukuthetha bucala kwindawo abangasayi kuphazanyiswa okanye baviwe. Ungakwazi ukusinceda sikwazi ukuqinisekisa ukuba uyakwazi ukuthetha kungekho namnye umntu omvayo, kuquka nawe. This is to ensure their privacy and confidentiality. Ndiyazi ukuba kunzima kwezinye zezindlu zethu, kodwa ungacinga ngexesha kunye nendawo apho oku kunokwenzeka? (jonga esi sithuba kunye nabo okanye mhlawumbi uyibeka njani into endaweni)"
Replace-All $d.Content "Thank you so much. We will call _________ at ____________ on the __________ . " "Enkosi kakhulu. We will call _________ at ____________ on the __________ . "

# --- Table cell (1,1): "If they haven't read it. Say..." block ---
# (scoped to the cell so it does not collide with the earlier
#  "If they answer no..." paragraph outside the table)
$tbl = $d.Tables.Item(1)
$cellA = $tbl.Cell(1,1).Range
Replace-All $cellA "If they " "Ukuba aba "
Replace-All $cellA "haven’t read" "khange"
Replace-All $cellA " it. Say…" " bayifunde. Ithi…"

# --- Table cell (1,2): "If they have read it. Say..." block ---
$cellB = $tbl.Cell(1,2).Range
Replace-All $cellB "If they " "Ukuba ba "
Replace-All $cellB "have read" "khe"
Replace-All $cellB " it. Say…" " bayifunda. Ithi…"
